$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Seating")
$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = "HT1"
$ws.Range("B1").Value = 7
$ws.Range("C1").Value = 8
$ws.Range("D1").Value = 13
$ws.Range("E1").Value = 10
$ws.Range("F1").Value = 9
$ws.Range("G1").Value = 11
$ws.Range("H1").Value = 12

$ws.Range("A1").Select()
